# Insert a new weekly price-observation row for "Membrillo" (Vega Modelo de
# Temuco) above the existing row 278, pushing the old rows 278-309 down to
# 279-310. The new row carries the same categorical attributes (market,
# region, product, quality, unit, origin) as the row that used to sit at 278,
# but with a newer date and updated volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 278:309 down to 279:310, duplicating formatting (incl. the date
# style on column D) from the row being pushed down - mirrors Excel's normal
# "Insert Cells/Row" behaviour.
$ws.Rows.Item(278).Insert()

# Populate the newly-created row 278 with the new observation.
$ws.Cells.Item(278, 1).Value  = 10                                      # A Mercado ID
$ws.Cells.Item(278, 2).Value  = "Vega Modelo de Temuco"                 # B Mercado
$ws.Cells.Item(278, 3).Value  = "La Araucanía"                          # C Región
$ws.Cells.Item(278, 4).Value  = 45132                                   # D Fecha
$ws.Cells.Item(278, 5).Value  = 9                                       # E Codreg
$ws.Cells.Item(278, 6).Value  = "Fruta"                                 # F Tipo
$ws.Cells.Item(278, 7).Value  = 100104                                  # G Producto ID
$ws.Cells.Item(278, 8).Value  = "Frutos de pepita"                      # H Producto
$ws.Cells.Item(278, 9).Value  = 100104003                               # I Categoría ID
$ws.Cells.Item(278, 10).Value = "Membrillo"                             # J Categoría
$ws.Cells.Item(278, 11).Value = "Champion"                              # K Variedad
$ws.Cells.Item(278, 12).Value = "Primera"                               # L Calidad
$ws.Cells.Item(278, 13).Value = 120                                     # M Volumen
$ws.Cells.Item(278, 14).Value = 16000                                   # N Precio mínimo
$ws.Cells.Item(278, 15).Value = 16000                                   # O Precio máximo
$ws.Cells.Item(278, 16).Value = 16000                                   # P Precio promedio ponderado
$ws.Cells.Item(278, 17).Value = "$/bandeja 18 kilos granel"             # Q Unidad de comercialización
$ws.Cells.Item(278, 18).Value = "Región de O'Higgins"                   # R Origen
$ws.Cells.Item(278, 19).Value = 889                                     # S Precio $/Kg
$ws.Cells.Item(278, 20).Value = 18                                      # T Kg / unidad
